$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

# Update header for the "Rolling ASE" column: Training Size 91 -> 95
$tbl.Cell(1, 4).Shape.TextFrame.TextRange.Text = "Rolling ASE (Training Size 95)"

# Update the Rolling ASE values (column 4) for each model row
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = "0.152"
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Text = "0.133"
$tbl.Cell(4, 4).Shape.TextFrame.TextRange.Text = "0.071"
$tbl.Cell(5, 4).Shape.TextFrame.TextRange.Text = "0.229"
$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "0.394"
